$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (A1) to the new header cells
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

$ws.Range("AC2:AC42").Value = 54
$ws.Range("AD2:AD42").Value = 61
$ws.Range("AE2:AE42").Value = 0
